# The presentation ships two theme parts:
#   ppt/theme/theme1.xml -> color scheme "Office"   (used by the Notes Master)
#   ppt/theme/theme2.xml -> color scheme "Integral" (used by the Slide Master / all slides)
#
# The authored change swaps the two color schemes between the parts, so the
# slides (and the rest of the deck that is actually driven by the Slide
# Master's theme) switch from the "Integral" palette to the "Office" palette.
#
# Drive this through the Slide Master's ColorScheme, which is the live,
# in-use theme for every slide/layout in the deck.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.ColorScheme

# Target palette ("Office" scheme, previously only used by the Notes Master)
# expressed as VBA-style BGR long values (R + G*256 + B*65536) for each of
# the twelve standard theme color slots, in their canonical order.
$scheme.Colors(1).RGB  = 0          # dk1      000000
$scheme.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$scheme.Colors(3).RGB  = 6968388    # dk2      44546A
$scheme.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$scheme.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$scheme.Colors(6).RGB  = 3243501    # accent2  ED7D31
$scheme.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$scheme.Colors(8).RGB  = 49407      # accent4  FFC000
$scheme.Colors(9).RGB  = 12874308   # accent5  4472C4
$scheme.Colors(10).RGB = 4697456    # accent6  70AD47
$scheme.Colors(11).RGB = 12673797   # hlink    0563C1
$scheme.Colors(12).RGB = 7491477    # folHlink 954F72
